# 40_Context3.pptx edit: remove the "40.3 Summary of Context 3" section-header
# slide (position 20) and renumber the trailing "40.x" headings that follow it.
#
# Before (positions 20-24):
#   20: [header]  "40.3 Summary of Context 3"
#   21: [content] "40.3 Summary of Context 3"
#   22: [header]  "40.4 Quiz"
#   23: [content] "40.4 Quiz"
#   24: [end]     "End of Chapter"
#
# After (positions 20-23):
#   20: [content] "40.2 Summary of Context 3"   (was pos 21)
#   21: [header]  "40.3 Quiz"                   (was pos 22)
#   22: [content] "40.3 Quiz"                   (was pos 23, title split into two runs)
#   23: [end]     "End of Chapter"               (was pos 24)

$p = $ppt.ActivePresentation

# 1. Delete the section-header slide "40.3 Summary of Context 3" at position 20.
#    Everything after it shifts up by one position.
$p.Slides.Item(20).Delete()

# 2. New slide 20 (was slide 21): "Summary of Context 3" content slide.
#    Its own title still reads "40.3 Summary of Context 3" -> bump to "40.2".
$slide20 = $p.Slides.Item(20)
$titleRange20 = $slide20.Shapes.Item(1).TextFrame.TextRange
$titleRange20.Runs(1, 1).Text = "40.2 Summary of Context 3"
# Cached slide-number field text: was "21", now "20".
$slide20.Shapes.Item(5).TextFrame.TextRange.Text = "20"

# 3. New slide 21 (was slide 22): "40.4 Quiz" section header -> "40.3 Quiz".
$slide21 = $p.Slides.Item(21)
$titleRange21 = $slide21.Shapes.Item(1).TextFrame.TextRange
$titleRange21.Runs(1, 1).Text = "40.3 Quiz"
# Cached slide-number field text stays logically "21" (matches new position).
$slide21.Shapes.Item(3).TextFrame.TextRange.Text = "21"

# 4. New slide 22 (was slide 23): "Quiz" content slide, title "40.4 Quiz" -> "40.3 Quiz"
#    split across two runs ("40.3 " + "Quiz"), matching the authored edit.
$slide22 = $p.Slides.Item(22)
$titleRange22 = $slide22.Shapes.Item(1).TextFrame.TextRange
$quizRun = $titleRange22.Runs(1, 1)
$quizRun.Text = "Quiz"
$quizRun.InsertBefore("40.3 ") | Out-Null
# Cached slide-number field text: was "23", now "22".
$slide22.Shapes.Item(5).TextFrame.TextRange.Text = "22"

# 5. New slide 23 (was slide 24): "End of Chapter" end slide, title unchanged.
$slide23 = $p.Slides.Item(23)
# Cached slide-number field text: was "24", now "23".
$slide23.Shapes.Item(3).TextFrame.TextRange.Text = "23"
